$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Julio de 2020 a las 21:58"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 4072509
$ws.Range("C4").Value = 43940
$ws.Range("D4").Value = 1909739
$ws.Range("E4").Value = 2017081
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 736
$ws.Range("H4").Value = 145689

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 1239684
$ws.Range("C6").Value = 45599
$ws.Range("D6").Value = 784266
$ws.Range("E6").Value = 425528
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 1120
$ws.Range("H6").Value = 29890

# Row 21: Alemania
$ws.Range("A21").Value = "Alemania"
$ws.Range("B21").Value = 204467
$ws.Range("C21").Value = 577
$ws.Range("D21").Value = 188600
$ws.Range("E21").Value = 6685
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 9182

# Row 31: Ecuador
$ws.Range("A31").Value = "Ecuador"
$ws.Range("B31").Value = 77257
$ws.Range("C31").Value = 1040
$ws.Range("D31").Value = 33125
$ws.Range("E31").Value = 38714
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 52
$ws.Range("H31").Value = 5418

# Row 56: Ghana
$ws.Range("A56").Value = "Ghana"
$ws.Range("B56").Value = 29672
$ws.Range("C56").Value = 683
$ws.Range("D56").Value = 26090
$ws.Range("E56").Value = 3429
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 153

# Row 77: Costa Rica
$ws.Range("A77").Value = "Costa Rica"
$ws.Range("B77").Value = 12361
$ws.Range("C77").Value = 550
$ws.Range("D77").Value = 3322
$ws.Range("E77").Value = 8968
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 71

# Row 81: Estado de Palestina
$ws.Range("A81").Value = "Estado de Palestina"
$ws.Range("B81").Value = 9398
$ws.Range("C81").Value = 170
$ws.Range("D81").Value = 1950
$ws.Range("E81").Value = 7382
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 66

# Row 92: Guayana Francesa
$ws.Range("A92").Value = "Guayana Francesa"
$ws.Range("B92").Value = 6883
$ws.Range("C92").Value = 32
$ws.Range("D92").Value = 5277
$ws.Range("E92").Value = 1567
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 39

# Row 150: Angola
$ws.Range("A150").Value = "Angola"
$ws.Range("B150").Value = 812
$ws.Range("C150").Value = 33
$ws.Range("D150").Value = 221
$ws.Range("E150").Value = 558
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 3
$ws.Range("H150").Value = 33

# Row 151: Jamaica
$ws.Range("A151").Value = "Jamaica"
$ws.Range("B151").Value = 810
$ws.Range("C151").Value = 1
$ws.Range("D151").Value = 709
$ws.Range("E151").Value = 91
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 10

# Row 152: Togo
$ws.Range("A152").Value = "Togo"
$ws.Range("B152").Value = 790
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 560
$ws.Range("E152").Value = 215
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 15

# Row 180: Gambia
$ws.Range("A180").Value = "Gambia"
$ws.Range("B180").Value = 146
$ws.Range("C180").Value = 34
$ws.Range("D180").Value = 57
$ws.Range("E180").Value = 84
$ws.Range("F180").Value = 0
$ws.Range("G180").Value = 1
$ws.Range("H180").Value = 5

# Row 181: Brunei
$ws.Range("A181").Value = "Brunei"
$ws.Range("B181").Value = 141
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 138
$ws.Range("E181").Value = 0
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 3

# Row 182: Trinidad yTobago
$ws.Range("A182").Value = "Trinidad yTobago"
$ws.Range("B182").Value = 139
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 124
$ws.Range("E182").Value = 7
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 8

# Row 199: Papua Nueva Guinea
$ws.Range("A199").Value = "Papua Nueva Guinea"
$ws.Range("B199").Value = 30
$ws.Range("C199").Value = 3
$ws.Range("D199").Value = 8
$ws.Range("E199").Value = 22
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0

# Row 200: Curazao
$ws.Range("A200").Value = "Curazao"
$ws.Range("B200").Value = 28
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 24
$ws.Range("E200").Value = 3
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 1

# Row 210: Islas Malvinas
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("B210").Value = 13
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 13
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Row 211: Groenlandia
$ws.Range("A211").Value = "Groenlandia"
$ws.Range("B211").Value = 13
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 13
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0
